# process confirm mail for 陈照明
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Candidates")

# Row 12 corresponds to candidate 陈照明 (Chen Zhaoming).
# L = 确认邮件 (confirmation email sent), M = 确认时间 (confirmation date),
# U = 确认意见 (confirmation note), V = 跟进意见 (follow-up note),
# W = 特定意见 (specific note)
$ws.Range("L12").Value = "Y"
$ws.Range("M12").Value = "'2019-04-28"
$ws.Range("U12").Value = "您投递的求职申请我们已经收到。我们将会在7日内完成对您简历的处理。感谢您对本职位的关注，我们将会尽快同您取得联系。"
$ws.Range("V12").Value = "尚未填写"
$ws.Range("W12").Value = "尚未填写"

$wb.Save()
